$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row before row 118, shifting existing rows 118:121 down to 119:122
$ws.Rows("118:118").Insert()

# Populate the newly inserted row with the new Fitch "BB-(EXP)" rating entry
$ws.Range("A118").Value = "Fitch"
$ws.Range("B118").Value = "BB-(EXP)"
$ws.Range("C118").Value = 9
$ws.Range("E118").Value = "Same as BB-"

# Update selection to reflect where the edit was made
$ws.Range("B119").Select()
